$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$searchRange = $ws.Range("G1:G259")

$firstAddress = $null
$found = $searchRange.Find($oldValue)
while ($found -ne $null) {
    if ($firstAddress -eq $null) {
        $firstAddress = $found.Address()
    } elseif ($found.Address() -eq $firstAddress) {
        break
    }
    $found.Value2 = $newValue
    $found = $searchRange.FindNext($found)
}
